$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Insert 5 new rows after row 8 (new Import statements for Spring
#    Expression Language classes). Use Insert + paste the formatting of
#    row 8 onto the newly inserted rows so the cell styles match (s=4 / s=6).
# ---------------------------------------------------------------------
$ws.Range("A9:D13").Insert(-4121) | Out-Null   # -4121 = xlShiftDown

$ws.Range("C8:D8").Copy() | Out-Null
$ws.Range("C9:D13").PasteSpecial(-4122) | Out-Null   # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C9").Value  = "Import"
$ws.Range("D9").Value  = "org.springframework.expression.EvaluationContext"
$ws.Range("C10").Value = "Import"
$ws.Range("D10").Value = "org.springframework.expression.Expression"
$ws.Range("C11").Value = "Import"
$ws.Range("D11").Value = "org.springframework.expression.ExpressionParser"
$ws.Range("C12").Value = "Import"
$ws.Range("D12").Value = "org.springframework.expression.spel.standard.SpelExpressionParser"
$ws.Range("C13").Value = "Import"
$ws.Range("D13").Value = "org.springframework.expression.spel.support.StandardEvaluationContext"

# ---------------------------------------------------------------------
# 2. Update the "Functions" body (now on row 14) to append the new
#    evalSpring() function. The appended block is formatted with its own
#    run (rPr) via Characters(), producing a 2-run rich string, matching
#    the target shared-string entry.
# ---------------------------------------------------------------------
$part1 = "function String dateFormat(String fmt)`n{`n  return LocalDate.now().toString(DateTimeFormat.forPattern(fmt));`n}`n`nfunction Date addDays(int days)`n{`n  Calendar cal = Calendar.getInstance();`n  cal.add(Calendar.DAY_OF_YEAR, days);`n  return cal.getTime();`n}`nfunction Boolean evalSpring(String expression, Object obj)`n"
$part2 = "{`n    ExpressionParser ep = new SpelExpressionParser();`n    Expression exp = ep.parseExpression(expression);`n    EvaluationContext ec = new StandardEvaluationContext();`n    Boolean evaluated = exp.getValue(ec, obj, Boolean.class);`n    return evaluated;`n}"

$ws.Range("D14").Value = ($part1 + $part2)

$startPos = $part1.Length + 1
$chars = $ws.Range("D14").Characters($startPos, $part2.Length)
$chars.Font.Size = 11
$chars.Font.Color = 0
$chars.Font.Name = "Calibri"

# ---------------------------------------------------------------------
# 3. Update the rule condition text that invokes evalSpring (row 20,
#    formerly row 15) and the Alfresco folder condition (row 23,
#    formerly row 18).
# ---------------------------------------------------------------------
$ws.Range("C20").Value = 'eval(evalSpring("$param", $caseFile))'
$ws.Range("C23").Value = 'container?.folder?.cmisFolderId == null'

# ---------------------------------------------------------------------
# 4. Restore the active selection recorded in the file.
# ---------------------------------------------------------------------
$ws.Range("C14").Select() | Out-Null
